$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "In kit?" header + column values (F) ---
$ws.Cells.Item(1,6).Value2 = "In kit?"
$ws.Cells.Item(1,6).Font.Bold = $true
$ws.Cells.Item(2,6).Value2 = "yes"
$ws.Cells.Item(3,6).Value2 = "yes"
$ws.Cells.Item(4,6).Value2 = "maybe"
$ws.Cells.Item(5,6).Value2 = "no"
$ws.Cells.Item(6,6).Value2 = "no"
$ws.Cells.Item(7,6).Value2 = "yes"
$ws.Cells.Item(8,6).Value2 = "no"

# --- Resistor rows 9-13 : now fully fleshed out component rows ---
$ws.Cells.Item(9,1).Value2 = "1k resistor"
$ws.Cells.Item(9,2).Value2 = "https://www.digikey.com/en/products/detail/yageo/CFR50SJT-52-1K/9099728"
$ws.Cells.Item(9,3).Value2 = 0.11
$ws.Cells.Item(9,4).Value2 = 7
$ws.Cells.Item(9,5).Value2 = 0.77
$ws.Cells.Item(9,6).Value2 = "yes"

$ws.Cells.Item(10,1).Value2 = "1M resistor"
$ws.Cells.Item(10,2).Value2 = "https://www.digikey.com/en/products/detail/yageo/CFR-25JT-52-1M/9098565"
$ws.Cells.Item(10,3).Value2 = 0.1
$ws.Cells.Item(10,4).Value2 = 1
$ws.Cells.Item(10,5).Value2 = 0.1
$ws.Cells.Item(10,6).Value2 = "yes"

$ws.Cells.Item(11,1).Value2 = "100k resistor"
$ws.Cells.Item(11,2).Value2 = "https://www.digikey.com/en/products/detail/yageo/CFR25SJT-52-100K/9098845"
$ws.Cells.Item(11,3).Value2 = 0.1
$ws.Cells.Item(11,4).Value2 = 4
$ws.Cells.Item(11,5).Value2 = 0.4
$ws.Cells.Item(11,6).Value2 = "yes"

$ws.Cells.Item(12,1).Value2 = "10k resistor"
$ws.Cells.Item(12,2).Value2 = "https://www.digikey.com/en/products/detail/yageo/CFR-12JT-52-10K/9098021"
$ws.Cells.Item(12,3).Value2 = 0.1
$ws.Cells.Item(12,4).Value2 = 1
$ws.Cells.Item(12,5).Value2 = 0.1
$ws.Cells.Item(12,6).Value2 = "yes"

$ws.Cells.Item(13,1).Value2 = "5k resistor"
$ws.Cells.Item(13,2).Value2 = "N/A"
$ws.Cells.Item(13,3).Value2 = "N/A"
$ws.Cells.Item(13,4).Value2 = 2
$ws.Cells.Item(13,5).Value2 = "N/A"
$ws.Cells.Item(13,6).Value2 = "yes"

# --- Pots rows 14-15 ---
$ws.Cells.Item(14,1).Value2 = "pots"
$ws.Cells.Item(14,2).Value2 = "https://www.digikey.com/en/products/detail/cui-devices/PTN10-B200SB20/16628511"
$ws.Cells.Item(14,3).Value2 = 0.65
$ws.Cells.Item(14,4).Value2 = 2
$ws.Cells.Item(14,5).Formula = "=C14*2"
$ws.Cells.Item(14,6).Value2 = "yes"

$ws.Cells.Item(15,2).Value2 = "https://www.digikey.com/en/products/detail/cui-devices/PTN10-E01SB20/16628490"
$ws.Cells.Item(15,3).Value2 = 0.65
$ws.Cells.Item(15,4).Value2 = 2
$ws.Cells.Item(15,5).Value2 = 1.3
$ws.Cells.Item(15,6).Value2 = "yes"

# --- Switch row 16 ---
$ws.Cells.Item(16,1).Value2 = "switch"
$ws.Cells.Item(16,2).Value2 = "https://www.digikey.com/en/products/detail/nidec-copal-electronics/ET310A12-Z/5086810"
$ws.Cells.Item(16,3).Value2 = 10.3
$ws.Cells.Item(16,4).Value2 = 1
$ws.Cells.Item(16,5).Value2 = 10.3
$ws.Cells.Item(16,6).Value2 = "no"

# --- NPN row 17 ---
$ws.Cells.Item(17,1).Value2 = "NPN"
$ws.Cells.Item(17,2).Value2 = "https://www.digikey.com/en/products/detail/stmicroelectronics/BUL138/1037756"
$ws.Cells.Item(17,3).Value2 = 1.59
$ws.Cells.Item(17,4).Value2 = 1
$ws.Cells.Item(17,5).Value2 = 1.59
$ws.Cells.Item(17,6).Value2 = "yes"

# --- PNP row 18 ---
$ws.Cells.Item(18,1).Value2 = "PNP"
$ws.Cells.Item(18,2).Value2 = "https://www.digikey.com/en/products/detail/onsemi/D45H8G/918452"
$ws.Cells.Item(18,3).Value2 = 0.9
$ws.Cells.Item(18,4).Value2 = 1
$ws.Cells.Item(18,5).Value2 = 0.9
$ws.Cells.Item(18,6).Value2 = "yes"

# --- New rows 19-20 : connectors ---
$ws.Cells.Item(19,1).Value2 = "8 pin conn"
$ws.Cells.Item(19,2).Value2 = "https://www.digikey.com/en/products/detail/adam-tech/ICS-308-T/9829299"
$ws.Cells.Item(19,3).Value2 = 0.1
$ws.Cells.Item(19,4).Value2 = 1
$ws.Cells.Item(19,5).Value2 = 0.1
$ws.Cells.Item(19,6).Value2 = "no"

$ws.Cells.Item(20,1).Value2 = "16 pin conn"
$ws.Cells.Item(20,2).Value2 = "https://www.digikey.com/en/products/detail/adam-tech/ICM-316-1-GT-HT/9833008"
$ws.Cells.Item(20,3).Value2 = 0.91
$ws.Cells.Item(20,4).Value2 = 2
$ws.Cells.Item(20,5).Value2 = 0.91
$ws.Cells.Item(20,6).Value2 = "yes"

# --- Apply currency number format to the newly populated C/E cells (rows 9-12, 14-20) ---
$ws.Range("C9:C12").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("E9:E12").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("C14:C20").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws.Range("E14:E20").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# --- Totals row 21 ---
$ws.Cells.Item(21,4).Value2 = "total cost:"
$ws.Cells.Item(21,5).Formula = "=SUM(E2:E20)"
$ws.Range("E21").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# --- Update the sheet selection to match the saved view ---
$ws.Range("B2:B20").Select()
